# semana 18 de 2025
# Update the Esperado (C), Observado (D) and valor p (E) columns on
# Sheet1 with the refreshed weekly Poisson figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Esperado (C), Observado (D), valor p (E).
# A $null placeholder means "leave this cell untouched".
$updates = @(
    @{ Row = 2;  C = 1;  D = 2;  E = 0.18 },
    @{ Row = 3;  C = $null; D = 1;  E = 0 },
    @{ Row = 4;  C = $null; D = 2;  E = 0.04 },
    @{ Row = 5;  C = $null; D = 7;  E = $null },
    @{ Row = 6;  C = 2;  D = $null; E = 0.09 },
    @{ Row = 8;  C = 37; D = 36; E = $null },
    @{ Row = 9;  C = $null; D = 0;  E = $null },
    @{ Row = 11; C = $null; D = 1;  E = 0.37 },
    @{ Row = 12; C = 3;  D = $null; E = 0.15 },
    @{ Row = 13; C = 8;  D = $null; E = $null },
    @{ Row = 14; C = 1;  D = $null; E = 0.37 },
    @{ Row = 16; C = 2;  D = $null; E = 0.14 },
    @{ Row = 17; C = $null; D = 4;  E = 0.01 },
    @{ Row = 19; C = 8;  D = 5;  E = 0.09 },
    @{ Row = 21; C = $null; D = 0;  E = 1 },
    @{ Row = 23; C = $null; D = 7;  E = 0.14 },
    @{ Row = 24; C = $null; D = 0;  E = 0.37 },
    @{ Row = 27; C = $null; D = 0;  E = $null },
    @{ Row = 28; C = $null; D = 1;  E = $null },
    @{ Row = 29; C = $null; D = 0;  E = 1 },
    @{ Row = 30; C = 1;  D = 2;  E = 0.18 },
    @{ Row = 32; C = $null; D = 7;  E = 0.14 },
    @{ Row = 33; C = 13; D = $null; E = 0.1 },
    @{ Row = 34; C = 7;  D = 9;  E = 0.1 }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.C) { $ws.Cells.Item($row, 3).Value = $u.C }
    if ($null -ne $u.D) { $ws.Cells.Item($row, 4).Value = $u.D }
    if ($null -ne $u.E) { $ws.Cells.Item($row, 5).Value = $u.E }
}
